$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename shared string label for row 9 (Thomas Hex -> Matthies Hex) and
# --- set new / shifted category labels for rows 4-31 (column B) ---
$ws.Range('B4').Value = 'Holden'
$ws.Range('B5').Value = 'Rizzie Spiral'
$ws.Range('B6').Value = 'RotRing OmegaMax-90'
$ws.Range('B7').Value = 'Equal Angle'
$ws.Range('B8').Value = 'Tilt Rotate'
$ws.Range('B9').Value = 'CLR'
$ws.Range('B10').Value = 'Rizzie Hex'
$ws.Range('B11').Value = 'Matthies Hex'
$ws.Range('B12').Value = 'Tilt Rotate_Partial'
$ws.Range('B13').Value = 'RotRing OmegaMax-60'
$ws.Range('B14').Value = 'Equal Angle_Partial'
$ws.Range('B15').Value = 'Rizzie Hex_Partial'
$ws.Range('B16').Value = 'ND Single'
$ws.Range('B17').Value = 'RD Single'
$ws.Range('B18').Value = 'TD Single'
$ws.Range('B19').Value = 'Morris Single'
$ws.Range('B20').Value = 'Ring Perpendicular to ND'
$ws.Range('B21').Value = 'Ring Perpendicular to RD'
$ws.Range('B22').Value = 'Ring Perpendicular to TD'
$ws.Range('B23').Value = 'OffsetFTD'
$ws.Range('B24').Value = 'OffsetATD'
$ws.Range('B25').Value = 'OffsetF45'
$ws.Range('B26').Value = 'OffsetA45'
$ws.Range('B27').Value = 'OffsetFRD'
$ws.Range('B28').Value = 'OffsetARD'
$ws.Range('B29').Value = 'Gaussian Quadrature'
$ws.Range('B30').Value = 'Michael-CCHex'
$ws.Range('B31').Value = 'Michael-SNHex'

# --- column A sequence numbers for new rows 30 and 31 ---
$ws.Range('A30').Value = 28
$ws.Range('A31').Value = 29

# --- numeric simulation data for rows 4-31, columns C:W ---
$ws.Range('C4').Value = 1.058601672565221
$ws.Range('D4').Value = 0.9846795121843712
$ws.Range('E4').Value = 0.9621234503712661
$ws.Range('F4').Value = 1.011291491402081
$ws.Range('G4').Value = 1.058601672565221
$ws.Range('H4').Value = 0.9857218425318898
$ws.Range('I4').Value = 1.037987128323161
$ws.Range('J4').Value = 0.9621234503712661
$ws.Range('K4').Value = 0.9621234503712661
$ws.Range('L4').Value = 0.9781772390827875
$ws.Range('M4').Value = 1.012551821969163
$ws.Range('N4').Value = 0.9621234503712661
$ws.Range('O4').Value = 0.9846795121843712
$ws.Range('P4').Value = 1.021640592374796
$ws.Range('Q4').Value = 0.9986156670767673
$ws.Range('R4').Value = 1.001801545040286
$ws.Range('S4').Value = 1.018611002239585
$ws.Range('T4').Value = 1.001801545040286
$ws.Range('U4').Value = 1.004489114272505
$ws.Range('V4').Value = 0.9960159814922577
$ws.Range('W4').Value = 1.003891769803743
$ws.Range('C5').Value = 1.142039601741266
$ws.Range('D5').Value = 0.9671335765114345
$ws.Range('E5').Value = 0.9006074138078776
$ws.Range('F5').Value = 1.030099745289099
$ws.Range('G5').Value = 1.142039601741266
$ws.Range('H5').Value = 0.9630215042172127
$ws.Range('I5').Value = 1.092812679342312
$ws.Range('J5').Value = 0.9006074138078776
$ws.Range('K5').Value = 0.9006074138078776
$ws.Range('L5').Value = 0.9475536382570774
$ws.Range('M5').Value = 1.03067978540747
$ws.Range('N5').Value = 0.9006074138078776
$ws.Range('O5').Value = 0.9671335765114345
$ws.Range('P5').Value = 1.05458658912635
$ws.Range('Q5').Value = 0.9989066809594522
$ws.Range('R5').Value = 1.003260197353526
$ws.Range('S5').Value = 1.04661765455339
$ws.Range('T5').Value = 1.003260197353526
$ws.Range('U5').Value = 1.010115094367012
$ws.Range('V5').Value = 0.988213558255185
$ws.Range('W5').Value = 1.009243493071719
$ws.Range('C6').Value = 1.059373575659167
$ws.Range('D6').Value = 0.9845877739158367
$ws.Range('E6').Value = 0.9614288697849082
$ws.Range('F6').Value = 1.011510662261824
$ws.Range('G6').Value = 1.059373575659167
$ws.Range('H6').Value = 0.9854726235174999
$ws.Range('I6').Value = 1.038506552962158
$ws.Range('J6').Value = 0.9614288697849082
$ws.Range('K6').Value = 0.9614288697849082
$ws.Range('L6').Value = 0.9779013419871293
$ws.Range('M6').Value = 1.012723767054321
$ws.Range('N6').Value = 0.9614288697849082
$ws.Range('O6').Value = 0.9845877739158367
$ws.Range('P6').Value = 1.021980674787502
$ws.Range('Q6').Value = 0.9986557704850787
$ws.Range('R6').Value = 1.001796739786637
$ws.Range('S6').Value = 1.018895038876441
$ws.Range('T6').Value = 1.001796739786637
$ws.Range('U6').Value = 1.004528496603558
$ws.Range('V6').Value = 0.9959085712398281
$ws.Range('W6').Value = 1.003938145892855
$ws.Range('C7').Value = 1.072493305641209
$ws.Range('D7').Value = 0.9811203140417862
$ws.Range('E7').Value = 0.9530157652377511
$ws.Range('F7').Value = 1.014014592752161
$ws.Range('G7').Value = 1.072493305641209
$ws.Range('H7').Value = 0.9822968729827088
$ws.Range('I7').Value = 1.047004608479827
$ws.Range('J7').Value = 0.9530157652377511
$ws.Range('K7').Value = 0.9530157652377511
$ws.Range('L7').Value = 0.9730117294020182
$ws.Range('M7').Value = 1.01553161507925
$ws.Range('N7').Value = 0.9530157652377511
$ws.Range('O7').Value = 0.9811203140417862
$ws.Range('P7').Value = 1.026806809841498
$ws.Range('Q7').Value = 0.998325964560518
$ws.Range('R7').Value = 1.002209794973582
$ws.Range('S7').Value = 1.023048411587415
$ws.Range('T7').Value = 1.002209794973582
$ws.Range('U7').Value = 1.005540249999999
$ws.Range('V7').Value = 0.9950353530475494
$ws.Range('W7').Value = 1.004811100452089
$ws.Range('C8').Value = 1.23963788681058
$ws.Range('D8').Value = 0.9361893109885945
$ws.Range('E8').Value = 0.8471768921643822
$ws.Range('F8').Value = 1.045430806587816
$ws.Range('G8').Value = 1.23963788681058
$ws.Range('H8').Value = 0.9422578673062515
$ws.Range('I8').Value = 1.15513847542685
$ws.Range('J8').Value = 0.8471768921643822
$ws.Range('K8').Value = 0.8471768921643822
$ws.Range('L8').Value = 0.910639016309732
$ws.Range('M8').Value = 1.051257976019466
$ws.Range('N8').Value = 0.8471768921643822
$ws.Range('O8').Value = 0.9361893109885945
$ws.Range('P8').Value = 1.087913598899587
$ws.Range('Q8').Value = 0.9937236435040304
$ws.Range('R8').Value = 1.007668029987852
$ws.Range('S8').Value = 1.075695057939547
$ws.Range('T8').Value = 1.007668029987852
$ws.Range('U8').Value = 1.018565516495756
$ws.Range('V8').Value = 0.9842877916294812
$ws.Range('W8').Value = 1.015966028951709
$ws.Range('C9').Value = 1.007693528795389
$ws.Range('D9').Value = 0.9981716536042885
$ws.Range('E9').Value = 0.9947020252352907
$ws.Range('F9').Value = 1.001599530313724
$ws.Range('G9').Value = 1.007693528795389
$ws.Range('H9').Value = 0.998023819780938
$ws.Range('I9').Value = 1.005018832996242
$ws.Range('J9').Value = 0.9947020252352907
$ws.Range('K9').Value = 0.9947020252352907
$ws.Range('L9').Value = 0.9971542047388436
$ws.Range('M9').Value = 1.001658865986806
$ws.Range('N9').Value = 0.9947020252352907
$ws.Range('O9').Value = 0.9981716536042885
$ws.Range('P9').Value = 1.002932591199839
$ws.Range('Q9').Value = 0.9999152597955472
$ws.Range('R9').Value = 1.000189069211656
$ws.Range('S9').Value = 1.002508016128828
$ws.Range('T9').Value = 1.000189069211656
$ws.Range('U9').Value = 1.000556518405444
$ws.Range('V9').Value = 0.999385619771413
$ws.Range('W9').Value = 1.00050280768144
$ws.Range('C10').Value = 1.00059128709084
$ws.Range('D10').Value = 0.9998503732997871
$ws.Range('E10').Value = 0.9996090141731878
$ws.Range('F10').Value = 1.000117103676048
$ws.Range('G10').Value = 1.00059128709084
$ws.Range('H10').Value = 0.9998531808991783
$ws.Range('I10').Value = 1.000384147178883
$ws.Range('J10').Value = 0.9996090141731878
$ws.Range('K10').Value = 0.9996090141731878
$ws.Range('L10').Value = 0.9997803288556475
$ws.Range('M10').Value = 1.000126944775042
$ws.Range('N10').Value = 0.9996090141731878
$ws.Range('O10').Value = 0.9998503732997871
$ws.Range('P10').Value = 1.000220830195313
$ws.Range('Q10').Value = 0.9999886590374145
$ws.Range('R10').Value = 1.000016891521271
$ws.Range('S10').Value = 1.000189535055223
$ws.Range('T10').Value = 1.000016891521271
$ws.Range('U10').Value = 1.000044404834714
$ws.Range('V10').Value = 0.9999573267024088
$ws.Range('W10').Value = 1.000039047493577
$ws.Range('C11').Value = 1.012901733186595
$ws.Range('D11').Value = 0.9969422340454663
$ws.Range('E11').Value = 0.9911007662758088
$ws.Range('F11').Value = 1.00268765540198
$ws.Range('G11').Value = 1.012901733186595
$ws.Range('H11').Value = 0.9966814192471165
$ws.Range('I11').Value = 1.008417812759083
$ws.Range('J11').Value = 0.9911007662758088
$ws.Range('K11').Value = 0.9911007662758088
$ws.Range('L11').Value = 0.995228589285313
$ws.Range('M11').Value = 1.00278235006737
$ws.Range('N11').Value = 0.9911007662758088
$ws.Range('O11').Value = 0.9969422340454663
$ws.Range('P11').Value = 1.004921983616031
$ws.Range('Q11').Value = 0.9998622920564182
$ws.Range('R11').Value = 1.00031491116929
$ws.Range('S11').Value = 1.004208772433144
$ws.Range('T11').Value = 1.00031491116929
$ws.Range('U11').Value = 1.00093177089381
$ws.Range('V11').Value = 0.9989655699702098
$ws.Range('W11').Value = 1.000842820033592
$ws.Range('C12').Value = 1.243247245765553
$ws.Range('D12').Value = 0.9351655290773299
$ws.Range('E12').Value = 0.8449865499734945
$ws.Range('F12').Value = 1.046074955786711
$ws.Range('G12').Value = 1.243247245765553
$ws.Range('H12').Value = 0.941422995690935
$ws.Range('I12').Value = 1.157464274623609
$ws.Range('J12').Value = 0.8449865499734945
$ws.Range('K12').Value = 0.8449865499734945
$ws.Range('L12').Value = 0.9092865068560856
$ws.Range('M12').Value = 1.052026242120317
$ws.Range('N12').Value = 0.8449865499734945
$ws.Range('O12').Value = 0.9351655290773299
$ws.Range('P12').Value = 1.089206387421441
$ws.Range('Q12').Value = 0.9935958855988236
$ws.Range('R12').Value = 1.007799774938793
$ws.Range('S12').Value = 1.0768130056544
$ws.Range('T12').Value = 1.007799774938793
$ws.Range('U12').Value = 1.018856391734174
$ws.Range('V12').Value = 0.984082423382038
$ws.Range('W12').Value = 1.016209287486754
$ws.Range('C13').Value = 1.036844430094802
$ws.Range('D13').Value = 0.9895618106513012
$ws.Range('E13').Value = 0.9776185198579166
$ws.Range('F13').Value = 1.006583553912415
$ws.Range('G13').Value = 1.036844430094802
$ws.Range('H13').Value = 0.991470602845962
$ws.Range('I13').Value = 1.023744000806459
$ws.Range('J13').Value = 0.9776185198579166
$ws.Range('K13').Value = 0.9776185198579166
$ws.Range('L13').Value = 0.9861948630512984
$ws.Range('M13').Value = 1.007843228877973
$ws.Range('N13').Value = 0.9776185198579166
$ws.Range('O13').Value = 0.9895618106513012
$ws.Range('P13').Value = 1.013203120373052
$ws.Range('Q13').Value = 0.9987025197646371
$ws.Range('R13').Value = 1.001341586868007
$ws.Range('S13').Value = 1.011416489874692
$ws.Range('T13').Value = 1.001341586868006
$ws.Range('U13').Value = 1.002966997370498
$ws.Range('V13').Value = 0.9978973018679816
$ws.Range('W13').Value = 1.002482626262266
$ws.Range('C14').Value = 1.073510688347368
$ws.Range('D14').Value = 0.9804578772210537
$ws.Range('E14').Value = 0.9530630057473659
$ws.Range('F14').Value = 1.013956891621054
$ws.Range('G14').Value = 1.073510688347368
$ws.Range('H14').Value = 0.9822692454947383
$ws.Range('I14').Value = 1.047595458221052
$ws.Range('J14').Value = 0.9530630057473659
$ws.Range('K14').Value = 0.9530630057473659
$ws.Range('L14').Value = 0.9725912490315769
$ws.Range('M14').Value = 1.01572569986316
$ws.Range('N14').Value = 0.9530630057473659
$ws.Range('O14').Value = 0.9804578772210537
$ws.Range('P14').Value = 1.026984282784211
$ws.Range('Q14').Value = 0.9980917885421069
$ws.Range('R14').Value = 1.002343857105263
$ws.Range('S14').Value = 1.023231421810527
$ws.Range('T14').Value = 1.002343857105263
$ws.Range('U14').Value = 1.005689317794737
$ws.Range('V14').Value = 0.9951640553852628
$ws.Range('W14').Value = 1.004896264443421
$ws.Range('C15').Value = 0.9526347222431559
$ws.Range('D15').Value = 1.011979104788094
$ws.Range('E15').Value = 1.031331895839029
$ws.Range('F15').Value = 0.99061513019074
$ws.Range('G15').Value = 0.9526347222431559
$ws.Range('H15').Value = 1.011764773345253
$ws.Range('I15').Value = 0.9692266824235301
$ws.Range('J15').Value = 1.031331895839029
$ws.Range('K15').Value = 1.031331895839029
$ws.Range('L15').Value = 1.017596035773496
$ws.Range('M15').Value = 0.9898306195768509
$ws.Range('N15').Value = 1.031331895839029
$ws.Range('O15').Value = 1.011979104788094
$ws.Range('P15').Value = 0.982306913515625
$ws.Range('Q15').Value = 1.000904862182473
$ws.Range('R15').Value = 0.998648574290093
$ws.Range('S15').Value = 0.9848148155360338
$ws.Range('T15').Value = 0.998648574290093
$ws.Range('U15').Value = 0.9964440856117824
$ws.Range('V15').Value = 1.003421647657232
$ws.Range('W15').Value = 0.9968723705225186
$ws.Range('C16').Value = 1.420376399999998
$ws.Range('D16').Value = 0.887797470000001
$ws.Range('E16').Value = 0.7323862499999992
$ws.Range('F16').Value = 1.079525899999999
$ws.Range('G16').Value = 1.420376399999998
$ws.Range('H16').Value = 0.8988550199999995
$ws.Range('I16').Value = 1.2721004
$ws.Range('J16').Value = 0.7323862499999992
$ws.Range('K16').Value = 0.7323862499999992
$ws.Range('L16').Value = 0.8432138100000013
$ws.Range('M16').Value = 1.0899016
$ws.Range('N16').Value = 0.7323862499999992
$ws.Range('O16').Value = 0.887797470000001
$ws.Range('P16').Value = 1.154086935
$ws.Range('Q16').Value = 0.9888495350000006
$ws.Range('R16').Value = 1.01352004
$ws.Range('S16').Value = 1.132691823333333
$ws.Range('T16').Value = 1.01352004
$ws.Range('U16').Value = 1.03261543
$ws.Range('V16').Value = 0.9725695939999996
$ws.Range('W16').Value = 1.02801960625
$ws.Range('C17').Value = 1.4203764
$ws.Range('D17').Value = 0.88779747
$ws.Range('E17').Value = 0.73238625
$ws.Range('F17').Value = 1.0795259
$ws.Range('G17').Value = 1.4203764
$ws.Range('H17').Value = 0.8988550199999998
$ws.Range('I17').Value = 1.2721004
$ws.Range('J17').Value = 0.73238625
$ws.Range('K17').Value = 0.73238625
$ws.Range('L17').Value = 0.8432138100000001
$ws.Range('M17').Value = 1.0899016
$ws.Range('N17').Value = 0.73238625
$ws.Range('O17').Value = 0.88779747
$ws.Range('P17').Value = 1.154086935
$ws.Range('Q17').Value = 0.9888495349999999
$ws.Range('R17').Value = 1.01352004
$ws.Range('S17').Value = 1.132691823333333
$ws.Range('T17').Value = 1.01352004
$ws.Range('U17').Value = 1.03261543
$ws.Range('V17').Value = 0.9725695940000001
$ws.Range('W17').Value = 1.02801960625
$ws.Range('C18').Value = 1.4203764
$ws.Range('D18').Value = 0.88779747
$ws.Range('E18').Value = 0.73238625
$ws.Range('F18').Value = 1.0795259
$ws.Range('G18').Value = 1.4203764
$ws.Range('H18').Value = 0.8988550199999998
$ws.Range('I18').Value = 1.2721004
$ws.Range('J18').Value = 0.73238625
$ws.Range('K18').Value = 0.73238625
$ws.Range('L18').Value = 0.8432138100000001
$ws.Range('M18').Value = 1.0899016
$ws.Range('N18').Value = 0.73238625
$ws.Range('O18').Value = 0.88779747
$ws.Range('P18').Value = 1.154086935
$ws.Range('Q18').Value = 0.9888495349999999
$ws.Range('R18').Value = 1.01352004
$ws.Range('S18').Value = 1.132691823333333
$ws.Range('T18').Value = 1.01352004
$ws.Range('U18').Value = 1.03261543
$ws.Range('V18').Value = 0.9725695940000001
$ws.Range('W18').Value = 1.02801960625
$ws.Range('C19').Value = 1.0209422
$ws.Range('D19').Value = 1.0016507
$ws.Range('E19').Value = 0.97379624
$ws.Range('F19').Value = 1.0085957
$ws.Range('G19').Value = 1.0209422
$ws.Range('H19').Value = 0.99093873
$ws.Range('I19').Value = 1.0148091
$ws.Range('J19').Value = 0.97379624
$ws.Range('K19').Value = 0.97379624
$ws.Range('L19').Value = 0.9929493
$ws.Range('M19').Value = 1.0049139
$ws.Range('N19').Value = 0.97379624
$ws.Range('O19').Value = 1.0016507
$ws.Range('P19').Value = 1.01129645
$ws.Range('Q19').Value = 1.0032823
$ws.Range('R19').Value = 0.9987963799999999
$ws.Range('S19').Value = 1.009168933333333
$ws.Range('T19').Value = 0.9987963799999999
$ws.Range('U19').Value = 1.00032576
$ws.Range('V19').Value = 0.9950198560000001
$ws.Range('W19').Value = 1.00107448375
$ws.Range('C20').Value = 1.157734744657534
$ws.Range('D20').Value = 0.9626598902739728
$ws.Range('E20').Value = 0.8911215853424657
$ws.Range('F20').Value = 1.032886836986301
$ws.Range('G20').Value = 1.157734744657534
$ws.Range('H20').Value = 0.9594032183561644
$ws.Range('I20').Value = 1.102922539726028
$ws.Range('J20').Value = 0.8911215853424657
$ws.Range('K20').Value = 0.8911215853424657
$ws.Range('L20').Value = 0.9416700165753423
$ws.Range('M20').Value = 1.034019244383562
$ws.Range('N20').Value = 0.8911215853424657
$ws.Range('O20').Value = 0.9626598902739728
$ws.Range('P20').Value = 1.060197317465754
$ws.Range('Q20').Value = 0.9983395673287676
$ws.Range('R20').Value = 1.003838740091324
$ws.Range('S20').Value = 1.051471293105023
$ws.Range('T20').Value = 1.003838740091324
$ws.Range('U20').Value = 1.011383866164384
$ws.Range('V20').Value = 0.9873314100000001
$ws.Range('W20').Value = 1.010302259537671
$ws.Range('C21').Value = 1.168102178421053
$ws.Range('D21').Value = 0.9597047947368422
$ws.Range('E21').Value = 0.8848557168421052
$ws.Range('F21').Value = 1.034727852631579
$ws.Range('G21').Value = 1.168102178421053
$ws.Range('H21').Value = 0.957013157894737
$ws.Range('I21').Value = 1.109600613157895
$ws.Range('J21').Value = 0.8848557168421052
$ws.Range('K21').Value = 0.8848557168421052
$ws.Range('L21').Value = 0.937783587368421
$ws.Range('M21').Value = 1.036225126842105
$ws.Range('N21').Value = 0.8848557168421052
$ws.Range('O21').Value = 0.9597047947368422
$ws.Range('P21').Value = 1.063903486578947
$ws.Range('Q21').Value = 0.9979649607894737
$ws.Range('R21').Value = 1.004220896666667
$ws.Range('S21').Value = 1.054677366666667
$ws.Range('T21').Value = 1.004220896666667
$ws.Range('U21').Value = 1.012221954210526
$ws.Range('V21').Value = 0.9867487067368422
$ws.Range('W21').Value = 1.011001628486842
$ws.Range('C22').Value = 1.168102178421053
$ws.Range('D22').Value = 0.9597047947368422
$ws.Range('E22').Value = 0.8848557168421052
$ws.Range('F22').Value = 1.034727852631579
$ws.Range('G22').Value = 1.168102178421053
$ws.Range('H22').Value = 0.957013157894737
$ws.Range('I22').Value = 1.109600613157895
$ws.Range('J22').Value = 0.8848557168421052
$ws.Range('K22').Value = 0.8848557168421052
$ws.Range('L22').Value = 0.937783587368421
$ws.Range('M22').Value = 1.036225126842105
$ws.Range('N22').Value = 0.8848557168421052
$ws.Range('O22').Value = 0.9597047947368422
$ws.Range('P22').Value = 1.063903486578947
$ws.Range('Q22').Value = 0.9979649607894737
$ws.Range('R22').Value = 1.004220896666667
$ws.Range('S22').Value = 1.054677366666667
$ws.Range('T22').Value = 1.004220896666667
$ws.Range('U22').Value = 1.012221954210526
$ws.Range('V22').Value = 0.9867487067368422
$ws.Range('W22').Value = 1.011001628486842
$ws.Range('C23').Value = 0.963735439783223
$ws.Range('D23').Value = 1.014369839509032
$ws.Range('E23').Value = 1.014747592515798
$ws.Range('F23').Value = 0.9961414537404596
$ws.Range('G23').Value = 0.963735439783223
$ws.Range('H23').Value = 1.006119641121082
$ws.Range('I23').Value = 0.9773388722374375
$ws.Range('J23').Value = 1.014747592515798
$ws.Range('K23').Value = 1.014747592515798
$ws.Range('L23').Value = 1.014017815243453
$ws.Range('M23').Value = 0.9925263924631762
$ws.Range('N23').Value = 1.014747592515798
$ws.Range('O23').Value = 1.014369839509032
$ws.Range('P23').Value = 0.9890526396461276
$ws.Range('Q23').Value = 1.003448115986104
$ws.Range('R23').Value = 0.9976176239360176
$ws.Range('S23').Value = 0.9902105572518104
$ws.Range('T23').Value = 0.9976176239360176
$ws.Range('U23').Value = 0.9963448160678072
$ws.Range('V23').Value = 1.000025371357405
$ws.Range('W23').Value = 0.9973746308267076
$ws.Range('C24').Value = 0.993306383595565
$ws.Range('D24').Value = 0.9994418808283211
$ws.Range('E24').Value = 1.008429562082902
$ws.Range('F24').Value = 0.9972331013818057
$ws.Range('G24').Value = 0.993306383595565
$ws.Range('H24').Value = 1.002913139413011
$ws.Range('I24').Value = 0.9952613996017942
$ws.Range('J24').Value = 1.008429562082902
$ws.Range('K24').Value = 1.008429562082902
$ws.Range('L24').Value = 1.002250371084115
$ws.Range('M24').Value = 0.9984275820655333
$ws.Range('N24').Value = 1.008429562082902
$ws.Range('O24').Value = 0.9994418808283211
$ws.Range('P24').Value = 0.996374132211943
$ws.Range('Q24').Value = 0.9989347314469272
$ws.Range('R24').Value = 1.000392608835596
$ws.Range('S24').Value = 0.997058615496473
$ws.Range('T24').Value = 1.000392608835596
$ws.Range('U24').Value = 0.9999013521430802
$ws.Range('V24').Value = 1.001606994131044
$ws.Range('W24').Value = 0.9996579275066309
$ws.Range('C25').Value = 0.92505585301876
$ws.Range('D25').Value = 1.016711760342585
$ws.Range('E25').Value = 1.053561470836624
$ws.Range('F25').Value = 0.9837156379439094
$ws.Range('G25').Value = 0.92505585301876
$ws.Range('H25').Value = 1.019860629219225
$ws.Range('I25').Value = 0.9509203947465752
$ws.Range('J25').Value = 1.053561470836624
$ws.Range('K25').Value = 1.053561470836624
$ws.Range('L25').Value = 1.027606115002572
$ws.Range('M25').Value = 0.9837746331294553
$ws.Range('N25').Value = 1.053561470836624
$ws.Range('O25').Value = 1.016711760342585
$ws.Range('P25').Value = 0.9708838066806723
$ws.Range('Q25').Value = 1.00024319673602
$ws.Range('R25').Value = 0.9984430280659896
$ws.Range('S25').Value = 0.9751807488302666
$ws.Range('T25').Value = 0.9984430280659896
$ws.Range('U25').Value = 0.9947759293318561
$ws.Range('V25').Value = 1.00653303763281
$ws.Range('W25').Value = 0.9951508117799632
$ws.Range('C26').Value = 0.9796151116138745
$ws.Range('D26').Value = 1.005927421620521
$ws.Range('E26').Value = 1.012112248171304
$ws.Range('F26').Value = 0.9964549868796057
$ws.Range('G26').Value = 0.9796151116138745
$ws.Range('H26').Value = 1.004634448438814
$ws.Range('I26').Value = 0.98688954647551
$ws.Range('J26').Value = 1.012112248171304
$ws.Range('K26').Value = 1.012112248171304
$ws.Range('L26').Value = 1.007653939649781
$ws.Range('M26').Value = 0.9956697340406946
$ws.Range('N26').Value = 1.012112248171304
$ws.Range('O26').Value = 1.005927421620521
$ws.Range('P26').Value = 0.9927712666171977
$ws.Range('Q26').Value = 1.000798577830608
$ws.Range('R26').Value = 0.9992182604685665
$ws.Range('S26').Value = 0.99373742242503
$ws.Range('T26').Value = 0.9992182604685665
$ws.Range('U26').Value = 0.9983311288615986
$ws.Range('V26').Value = 1.00108735272354
$ws.Range('W26').Value = 0.9986196796112632
$ws.Range('C27').Value = 0.963735439783223
$ws.Range('D27').Value = 1.014369839509032
$ws.Range('E27').Value = 1.014747592515798
$ws.Range('F27').Value = 0.9961414537404596
$ws.Range('G27').Value = 0.963735439783223
$ws.Range('H27').Value = 1.006119641121082
$ws.Range('I27').Value = 0.9773388722374375
$ws.Range('J27').Value = 1.014747592515798
$ws.Range('K27').Value = 1.014747592515798
$ws.Range('L27').Value = 1.014017815243453
$ws.Range('M27').Value = 0.9925263924631761
$ws.Range('N27').Value = 1.014747592515798
$ws.Range('O27').Value = 1.014369839509032
$ws.Range('P27').Value = 0.9890526396461276
$ws.Range('Q27').Value = 1.003448115986104
$ws.Range('R27').Value = 0.9976176239360176
$ws.Range('S27').Value = 0.9902105572518104
$ws.Range('T27').Value = 0.9976176239360176
$ws.Range('U27').Value = 0.9963448160678072
$ws.Range('V27').Value = 1.000025371357405
$ws.Range('W27').Value = 0.9973746308267076
$ws.Range('C28').Value = 0.9933063835955649
$ws.Range('D28').Value = 0.9994418808283211
$ws.Range('E28').Value = 1.008429562082902
$ws.Range('F28').Value = 0.9972331013818057
$ws.Range('G28').Value = 0.9933063835955649
$ws.Range('H28').Value = 1.002913139413011
$ws.Range('I28').Value = 0.9952613996017942
$ws.Range('J28').Value = 1.008429562082902
$ws.Range('K28').Value = 1.008429562082902
$ws.Range('L28').Value = 1.002250371084115
$ws.Range('M28').Value = 0.9984275820655333
$ws.Range('N28').Value = 1.008429562082902
$ws.Range('O28').Value = 0.9994418808283211
$ws.Range('P28').Value = 0.996374132211943
$ws.Range('Q28').Value = 0.9989347314469272
$ws.Range('R28').Value = 1.000392608835596
$ws.Range('S28').Value = 0.997058615496473
$ws.Range('T28').Value = 1.000392608835596
$ws.Range('U28').Value = 0.9999013521430802
$ws.Range('V28').Value = 1.001606994131044
$ws.Range('W28').Value = 0.9996579275066309
$ws.Range('C29').Value = 1.024628540296642
$ws.Range('D29').Value = 0.9940048457138401
$ws.Range('E29').Value = 0.9832929756841889
$ws.Range('F29').Value = 1.005029388252175
$ws.Range('G29').Value = 1.024628540296642
$ws.Range('H29').Value = 0.993752856513907
$ws.Range('I29').Value = 1.016041686398808
$ws.Range('J29').Value = 0.9832929756841889
$ws.Range('K29').Value = 0.9832929756841889
$ws.Range('L29').Value = 0.9908751033011152
$ws.Range('M29').Value = 1.0053018174538
$ws.Range('N29').Value = 0.9832929756841889
$ws.Range('O29').Value = 0.9940048457138401
$ws.Range('P29').Value = 1.009316693005241
$ws.Range('Q29').Value = 0.9996533315838203
$ws.Range('R29').Value = 1.00064212056489
$ws.Range('S29').Value = 1.007978401154761
$ws.Range('T29').Value = 1.00064212056489
$ws.Range('U29').Value = 1.001807044787118
$ws.Range('V29').Value = 0.9981042309665321
$ws.Range('W29').Value = 1.00161590170181
$ws.Range('C30').Value = 1.01068436835141
$ws.Range('D30').Value = 0.9990660449804287
$ws.Range('E30').Value = 0.9897888390056322
$ws.Range('F30').Value = 1.003248647264372
$ws.Range('G30').Value = 1.01068436835141
$ws.Range('H30').Value = 0.9963638310877883
$ws.Range('I30').Value = 1.007247819337628
$ws.Range('J30').Value = 0.9897888390056322
$ws.Range('K30').Value = 0.9897888390056322
$ws.Range('L30').Value = 0.9962164034494032
$ws.Range('M30').Value = 1.002400226879678
$ws.Range('N30').Value = 0.9897888390056322
$ws.Range('O30').Value = 0.9990660449804287
$ws.Range('P30').Value = 1.004875206665919
$ws.Range('Q30').Value = 1.000733135930054
$ws.Range('R30').Value = 0.9998464174458235
$ws.Range('S30').Value = 1.004050213403839
$ws.Range('T30').Value = 0.9998464174458236
$ws.Range('U30').Value = 1.000484869804287
$ws.Range('V30').Value = 0.9983456636445563
$ws.Range('W30').Value = 1.000627022544543
$ws.Range('C31').Value = 0.93427825215649
$ws.Range('D31').Value = 1.019018088069463
$ws.Range('E31').Value = 1.039214206429792
$ws.Range('F31').Value = 0.9885117662537631
$ws.Range('G31').Value = 0.93427825215649
$ws.Range('H31').Value = 1.014992848993604
$ws.Range('I31').Value = 0.9577153773118796
$ws.Range('J31').Value = 1.039214206429792
$ws.Range('K31').Value = 1.039214206429792
$ws.Range('L31').Value = 1.024666954473673
$ws.Range('M31').Value = 0.9860335068299773
$ws.Range('N31').Value = 1.039214206429792
$ws.Range('O31').Value = 1.019018088069463
$ws.Range('P31').Value = 0.9766481701129768
$ws.Range('Q31').Value = 1.00252579744972
$ws.Range('R31').Value = 0.9975035155519151
$ws.Range('S31').Value = 0.9797766156853104
$ws.Range('T31').Value = 0.997503515551915
$ws.Range('U31').Value = 0.9946360133714305
$ws.Range('V31').Value = 1.003551651983103
$ws.Range('W31').Value = 0.9955538750648303
